$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.196078431372549
$ws.Range("C2").Value = 0.5450980392156862
$ws.Range("J2").Value = 0.007843137254901961
$ws.Range("P2").Value = 0.1686274509803922
$ws.Range("S2").Value = 0.08235294117647059
$ws.Range("C3").Value = 0.02068965517241379
$ws.Range("J3").Value = 0.02758620689655172
$ws.Range("P3").Value = 0.6896551724137931
$ws.Range("S3").Value = 0.2620689655172414
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("O4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.2121212121212121
$ws.Range("B6").Value = 0.04444444444444445
$ws.Range("D6").Value = 0.02222222222222222
$ws.Range("F6").Value = 0.03333333333333333
$ws.Range("J6").Value = 0.2777777777777778
$ws.Range("O6").Value = 0.01111111111111111
$ws.Range("Q6").Value = 0.1333333333333333
$ws.Range("R6").Value = 0.06111111111111111
$ws.Range("S6").Value = 0.4166666666666667
$ws.Range("B7").Value = 0.09090909090909091
$ws.Range("D7").Value = 0.02392344497607655
$ws.Range("F7").Value = 0.01435406698564593
$ws.Range("J7").Value = 0.1770334928229665
$ws.Range("O7").Value = 0.02392344497607655
$ws.Range("Q7").Value = 0.1339712918660287
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.4449760765550239
$ws.Range("B8").Value = 0.08430913348946135
$ws.Range("D8").Value = 0.00936768149882904
$ws.Range("F8").Value = 0.04215456674473068
$ws.Range("J8").Value = 0.1358313817330211
$ws.Range("Q8").Value = 0.17096018735363
$ws.Range("R8").Value = 0.07494145199063232
$ws.Range("S8").Value = 0.4824355971896955
$ws.Range("B9").Value = 0.07344632768361582
$ws.Range("D9").Value = 0.01694915254237288
$ws.Range("F9").Value = 0.0847457627118644
$ws.Range("J9").Value = 0.1129943502824859
$ws.Range("O9").Value = 0.01129943502824859
$ws.Range("Q9").Value = 0.1807909604519774
$ws.Range("R9").Value = 0.1016949152542373
$ws.Range("S9").Value = 0.4180790960451977
$ws.Range("B10").Value = 0.1082474226804124
$ws.Range("D10").Value = 0.01632302405498282
$ws.Range("E10").Value = 0.000859106529209622
$ws.Range("F10").Value = 0.07302405498281787
$ws.Range("J10").Value = 0.1202749140893471
$ws.Range("O10").Value = 0.01374570446735395
$ws.Range("Q10").Value = 0.1829896907216495
$ws.Range("R10").Value = 0.0781786941580756
$ws.Range("S10").Value = 0.4063573883161512
$ws.Range("G11").Value = 0.1512605042016807
$ws.Range("J11").Value = 0.08403361344537816
$ws.Range("K11").Value = 0.2128851540616246
$ws.Range("L11").Value = 0.5378151260504201
$ws.Range("S11").Value = 0.01400560224089636
$ws.Range("G12").Value = 0.6881188118811881
$ws.Range("J12").Value = 0.2178217821782178
$ws.Range("K12").Value = 0.009900990099009901
$ws.Range("L12").Value = 0.0396039603960396
$ws.Range("S12").Value = 0.04455445544554455
$ws.Range("G13").Value = 0.6363636363636364
$ws.Range("J13").Value = 0.3636363636363636
$ws.Range("F15").Value = 0.01507537688442211
$ws.Range("H15").Value = 0.1457286432160804
$ws.Range("I15").Value = 0.09045226130653267
$ws.Range("J15").Value = 0.3718592964824121
$ws.Range("K15").Value = 0.07537688442211055
$ws.Range("M15").Value = 0.01507537688442211
$ws.Range("O15").Value = 0.04522613065326633
$ws.Range("S15").Value = 0.2412060301507538
$ws.Range("H16").Value = 0.191358024691358
$ws.Range("I16").Value = 0.07407407407407407
$ws.Range("J16").Value = 0.382716049382716
$ws.Range("K16").Value = 0.08641975308641975
$ws.Range("M16").Value = 0.006172839506172839
$ws.Range("O16").Value = 0.03703703703703703
$ws.Range("S16").Value = 0.2222222222222222
$ws.Range("F17").Value = 0.01351351351351351
$ws.Range("H17").Value = 0.1621621621621622
$ws.Range("I17").Value = 0.08918918918918919
$ws.Range("J17").Value = 0.3567567567567568
$ws.Range("K17").Value = 0.1324324324324324
$ws.Range("M17").Value = 0.02162162162162162
$ws.Range("O17").Value = 0.1
$ws.Range("S17").Value = 0.1243243243243243
$ws.Range("F18").Value = 0.01183431952662722
$ws.Range("H18").Value = 0.2189349112426036
$ws.Range("I18").Value = 0.05917159763313609
$ws.Range("J18").Value = 0.4260355029585799
$ws.Range("K18").Value = 0.07692307692307693
$ws.Range("M18").Value = 0.005917159763313609
$ws.Range("N18").Value = 0.005917159763313609
$ws.Range("O18").Value = 0.05917159763313609
$ws.Range("S18").Value = 0.136094674556213
$ws.Range("F19").Value = 0.01607963246554365
$ws.Range("H19").Value = 0.2082695252679939
$ws.Range("I19").Value = 0.08116385911179173
$ws.Range("J19").Value = 0.3415007656967841
$ws.Range("K19").Value = 0.1424196018376723
$ws.Range("M19").Value = 0.01607963246554365
$ws.Range("N19").Value = 0.003062787136294028
$ws.Range("O19").Value = 0.06431852986217458
$ws.Range("S19").Value = 0.1271056661562021
